# 31May2022 Selenium DataDriven Part 5
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Fill in the "Result" column (S) for the data rows that were still blank.
$ws.Range("S5").Value = "Pass"
$ws.Range("S8").Value = "Pass"

# Update the selection on the sheet to match the new working range.
$ws.Range("S2:S7").Select()

# Collapse column I (was expanded in the source workbook).
$ws.Columns.Item(9).Collapsed = $true
